# Daily attendance processing - 2026-01-07 06:45:59
# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (G) for every row that currently has that value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$firstRow = $used.Row

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cell = $ws.Range("G$r")
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
